# Auto-generated Excel COM-interop script applying scheduled-runner updates
# to the Ixion_Profits workbook (profit recalculation across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H129").Value = 1187.069
$ws.Range("I129").Value = 832.25
$ws.Range("K129").Value = 2496.75
$ws.Range("M129").Value = 2503.25

$ws.Range("H132").Value = 1915.9259
$ws.Range("I132").Value = 1672.1923
$ws.Range("J132").Value = 8253
$ws.Range("K132").Value = 5016.5769
$ws.Range("L132").Value = 24759
$ws.Range("M132").Value = -2486.5769
$ws.Range("N132").Value = -29819

$ws.Range("H136").Value = 47425
$ws.Range("J136").Value = 47425
$ws.Range("L136").Value = 47425
$ws.Range("N136").Value = -57625

$ws.Range("H139").Value = 47314.285
$ws.Range("J139").Value = 47314.285
$ws.Range("L139").Value = 47314.285
$ws.Range("N139").Value = -57594.285

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H34").Value = 11900
$ws.Range("J34").Value = 11900
$ws.Range("L34").Value = 11900
$ws.Range("N34").Value = -12442

$ws.Range("H45").Value = 234668
$ws.Range("I45").Value = 420502.4
$ws.Range("K45").Value = 420502.4
$ws.Range("M45").Value = -420125.4

$ws.Range("H132").Value = 1518472.9
$ws.Range("I132").Value = 1793.579
$ws.Range("J132").Value = 3576823.2
$ws.Range("K132").Value = 5380.737
$ws.Range("L132").Value = 10730469.6
$ws.Range("M132").Value = -2850.737
$ws.Range("N132").Value = -10735529.6

$ws.Range("H135").Value = 56422.8
$ws.Range("J135").Value = 56422.8
$ws.Range("L135").Value = 56422.8
$ws.Range("N135").Value = -66562.8

$ws.Range("H138").Value = 42943
$ws.Range("J138").Value = 42943
$ws.Range("L138").Value = 42943
$ws.Range("N138").Value = -53223

$ws.Range("H139").Value = 66150
$ws.Range("J139").Value = 66150
$ws.Range("L139").Value = 66150
$ws.Range("N139").Value = -76430

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H11").Value = 581.8182
$ws.Range("I11").Value = 540
$ws.Range("K11").Value = 540
$ws.Range("M11").Value = -400

$ws.Range("H81").Value = 39005.715
$ws.Range("J81").Value = 39005.715
$ws.Range("L81").Value = 39005.715
$ws.Range("N81").Value = -41127.715

$ws.Range("H84").Value = 39005.715
$ws.Range("J84").Value = 39005.715
$ws.Range("L84").Value = 117017.145
$ws.Range("N84").Value = -127625.145

$ws.Range("H94").Value = 1956.1875
$ws.Range("I94").Value = 1545.4445
$ws.Range("J94").Value = 2484.2856
$ws.Range("K94").Value = 1545.4445
$ws.Range("L94").Value = 2484.2856
$ws.Range("M94").Value = -1094.4445
$ws.Range("N94").Value = -3386.2856

$ws.Range("H137").Value = 66630.336
$ws.Range("J137").Value = 57956.4
$ws.Range("L137").Value = 57956.4
$ws.Range("N137").Value = -68156.39999999999

$ws.Range("H138").Value = 60692
$ws.Range("J138").Value = 60692
$ws.Range("L138").Value = 60692
$ws.Range("N138").Value = -70972

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H16").Value = 4275646
$ws.Range("I16").Value = 6994726.5
$ws.Range("J16").Value = 2804.8572
$ws.Range("K16").Value = 6994726.5
$ws.Range("L16").Value = 2804.8572
$ws.Range("M16").Value = -6994439.5
$ws.Range("N16").Value = -3378.8572

$ws.Range("H31").Value = 10282.8125
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 10282.8125
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 10282.8125
$ws.Range("N31").Value = -10872.8125
$ws.Range("M31").ClearContents()

$ws.Range("H34").Value = 10282.8125
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 10282.8125
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 10282.8125
$ws.Range("N34").Value = -10686.8125
$ws.Range("M34").ClearContents()

$ws.Range("H113").Value = 4275646
$ws.Range("I113").Value = 6994726.5
$ws.Range("J113").Value = 2804.8572
$ws.Range("K113").Value = 6994726.5
$ws.Range("L113").Value = 2804.8572
$ws.Range("M113").Value = -6992556.5
$ws.Range("N113").Value = -7144.8572

$ws.Range("H138").Value = 66500
$ws.Range("J138").Value = 66500
$ws.Range("L138").Value = 66500
$ws.Range("N138").Value = -76780

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H123").Value = 6676.6665
$ws.Range("I123").Value = 1030
$ws.Range("J123").Value = 9500
$ws.Range("K123").Value = 3090
$ws.Range("L123").Value = 28500
$ws.Range("M123").Value = -640
$ws.Range("N123").Value = -33400

$ws.Range("H131").Value = 1516220.5
$ws.Range("I131").Value = 4545970
$ws.Range("J131").Value = 1345.6364
$ws.Range("K131").Value = 13637910
$ws.Range("L131").Value = 4036.9092
$ws.Range("M131").Value = -13632870
$ws.Range("N131").Value = -14116.9092

$ws.Range("H132").Value = 2631.6843
$ws.Range("I132").Value = 1255.2222
$ws.Range("K132").Value = 11296.9998
$ws.Range("M132").Value = -8766.9998

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H20").Value = 5006316
$ws.Range("I20").Value = 20000000
$ws.Range("K20").Value = 20000000
$ws.Range("M20").Value = -19999755

$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()

$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").ClearContents()

$ws.Range("H113").Value = 76925130
$ws.Range("I113").Value = 142859420
$ws.Range("K113").Value = 142859420
$ws.Range("M113").Value = -142857250

$ws.Range("H119").Value = 32000
$ws.Range("J119").Value = 32000
$ws.Range("L119").Value = 32000
$ws.Range("N119").Value = -41676

$ws.Range("H132").Value = 10628.211
$ws.Range("I132").Value = 8819.823
$ws.Range("J132").Value = 25999.5
$ws.Range("K132").Value = 26459.469
$ws.Range("L132").Value = 77998.5
$ws.Range("M132").Value = -23929.469
$ws.Range("N132").Value = -83058.5

$ws.Range("H136").Value = 60442
$ws.Range("J136").Value = 60442
$ws.Range("L136").Value = 181326
$ws.Range("N136").Value = -186426

$ws.Range("H138").Value = 50429
$ws.Range("J138").Value = 50429
$ws.Range("L138").Value = 50429
$ws.Range("N138").Value = -60709

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H25").Value = 6001.1665
$ws.Range("I25").Value = 1201.4
$ws.Range("K25").Value = 1201.4
$ws.Range("M25").Value = -971.4000000000001

$ws.Range("H122").Value = 8150136
$ws.Range("I122").Value = 7944596
$ws.Range("K122").Value = 23833788
$ws.Range("M122").Value = -23831338

$ws.Range("H136").Value = 5844.4146
$ws.Range("I136").Value = 1710.5186
$ws.Range("J136").Value = 13816.929
$ws.Range("K136").Value = 5131.5558
$ws.Range("L136").Value = 41450.787
$ws.Range("M136").Value = -2581.5558
$ws.Range("N136").Value = -46550.787

$ws.Range("H138").Value = 75365.42999999999
$ws.Range("J138").Value = 75365.42999999999
$ws.Range("L138").Value = 75365.42999999999
$ws.Range("N138").Value = -85645.42999999999

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

$ws.Range("H18").Value = 500003620
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 500003620
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 500003620
$ws.Range("N18").Value = -500003966
$ws.Range("M18").ClearContents()

$ws.Range("H20").Value = 5003605
$ws.Range("J20").Value = 7210
$ws.Range("L20").Value = 7210
$ws.Range("N20").Value = -7690

$ws.Range("H43").Value = 9666.666999999999
$ws.Range("J43").Value = 9666.666999999999
$ws.Range("L43").Value = 9666.666999999999
$ws.Range("N43").Value = -9964.666999999999

$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H132").Value = 1400.7142
$ws.Range("I132").Value = 556.5925999999999
$ws.Range("J132").Value = 4249.625
$ws.Range("K132").Value = 1669.7778
$ws.Range("L132").Value = 12748.875
$ws.Range("M132").Value = 860.2222000000002
$ws.Range("N132").Value = -17808.875

$ws.Range("H137").Value = 45678.75
$ws.Range("J137").Value = 45678.75
$ws.Range("L137").Value = 45678.75
$ws.Range("N137").Value = -55878.75
